$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from column P (last existing header) into new header columns Q:T
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:T1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Header row (row 1): new "Estación más cercana 6/7" and "Inicio..." columns
$ws.Range("Q1").Value = "Estación más cercana 6"
$ws.Range("R1").Value = "Estación más cercana 7"
$ws.Range("S1").Value = "Inicio estación más cercana 6"
$ws.Range("T1").Value = "Inicio estación más cercana 7"

# Data rows 2-40: Estación más cercana 6/7 + Inicio estación más cercana 6/7
# Row 2
$ws.Range("Q2").Value = "PX03"
$ws.Range("R2").Value = "PATCX"
$ws.Range("S2").Value = "2014-01-20T23:32:58"
$ws.Range("T2").Value = "2014-01-20T23:32:59"
# Row 3
$ws.Range("Q3").Value = "TLL"
$ws.Range("R3").Value = "CO03"
$ws.Range("S3").Value = "2014-03-14T15:31:11"
$ws.Range("T3").Value = "2014-03-14T15:31:13"
# Row 4
$ws.Range("Q4").Value = "PATCX"
$ws.Range("R4").Value = "PB11"
$ws.Range("S4").Value = "2014-03-17T10:05:41"
$ws.Range("T4").Value = "2014-03-17T10:05:41"
# Row 5
$ws.Range("Q5").Value = "PB04"
$ws.Range("R5").Value = "PB10"
$ws.Range("S5").Value = "2014-03-20T11:40:28"
$ws.Range("T5").Value = "2014-03-20T11:40:29"
# Row 6
$ws.Range("Q6").Value = "PSGCX"
$ws.Range("R6").Value = "HMBCX"
$ws.Range("S6").Value = "2014-03-23T07:03:20"
$ws.Range("T6").Value = "2014-03-23T07:03:22"
# Row 7
$ws.Range("Q7").Value = "PB12"
$ws.Range("R7").Value = "TA01"
$ws.Range("S7").Value = "2014-03-24T17:59:23"
$ws.Range("T7").Value = "2014-03-24T17:59:24"
# Row 8
$ws.Range("Q8").Value = "TA01"
$ws.Range("R8").Value = "HMBCX"
$ws.Range("S8").Value = "2014-03-26T02:14:57"
$ws.Range("T8").Value = "2014-03-26T02:14:57"
# Row 9
$ws.Range("Q9").Value = "PB11"
$ws.Range("R9").Value = "PX02"
$ws.Range("S9").Value = "2014-03-29T18:00:41"
$ws.Range("T9").Value = "2014-03-29T18:00:41"
# Row 10
$ws.Range("Q10").Value = "PSGCX"
$ws.Range("R10").Value = "PX03"
$ws.Range("S10").Value = "2014-04-03T19:57:05"
$ws.Range("T10").Value = "2014-04-03T19:57:05"
# Row 11
$ws.Range("Q11").Value = "PB12"
$ws.Range("R11").Value = "TA01"
$ws.Range("S11").Value = "2014-04-03T21:32:33"
$ws.Range("T11").Value = "2014-04-03T21:32:34"
# Row 12
$ws.Range("Q12").Value = "PATCX"
$ws.Range("R12").Value = "PX02"
$ws.Range("S12").Value = "2014-04-11T19:08:47"
$ws.Range("T12").Value = "2014-04-11T19:08:47"
# Row 13
$ws.Range("Q13").Value = "PB11"
$ws.Range("R13").Value = "TA01"
$ws.Range("S13").Value = "2014-04-15T01:59:46"
$ws.Range("T13").Value = "2014-04-15T01:59:47"
# Row 14
$ws.Range("Q14").Value = "PX02"
$ws.Range("R14").Value = "PB11"
$ws.Range("S14").Value = "2014-04-21T07:33:55"
$ws.Range("T14").Value = "2014-04-21T07:33:56"
# Row 15
$ws.Range("Q15").Value = "PB11"
$ws.Range("R15").Value = "PX03"
$ws.Range("S15").Value = "2014-04-25T13:41:01"
$ws.Range("T15").Value = "2014-04-25T13:41:02"
# Row 16
$ws.Range("Q16").Value = "PB11"
$ws.Range("R16").Value = "PB08"
$ws.Range("S16").Value = "2014-04-27T13:23:21"
$ws.Range("T16").Value = "2014-04-27T13:23:23"
# Row 17
$ws.Range("Q17").Value = "PATCX"
$ws.Range("R17").Value = "PB11"
$ws.Range("S17").Value = "2014-05-19T06:03:08"
$ws.Range("T17").Value = "2014-05-19T06:03:09"
# Row 18
$ws.Range("Q18").Value = "PB11"
$ws.Range("R18").Value = "PB08"
$ws.Range("S18").Value = "2014-05-24T08:03:14"
$ws.Range("T18").Value = "2014-05-24T08:03:16"
# Row 19
$ws.Range("Q19").Value = "PATCX"
$ws.Range("R19").Value = "PB01"
$ws.Range("S19").Value = "2014-05-28T00:41:35"
$ws.Range("T19").Value = "2014-05-28T00:41:37"
# Row 20
$ws.Range("Q20").Value = "PX03"
$ws.Range("R20").Value = "PB11"
$ws.Range("S20").Value = "2014-06-07T07:48:45"
$ws.Range("T20").Value = "2014-06-07T07:48:47"
# Row 21
$ws.Range("Q21").Value = "PX05"
$ws.Range("R21").Value = "PB20"
$ws.Range("S21").Value = "2014-06-20T06:29:43"
$ws.Range("T21").Value = "2014-06-20T06:29:46"
# Row 22
$ws.Range("Q22").Value = "PX03"
$ws.Range("R22").Value = "PSGCX"
$ws.Range("S22").Value = "2014-07-02T12:09:04"
$ws.Range("T22").Value = "2014-07-02T12:09:05"
# Row 23
$ws.Range("Q23").Value = "VA01"
$ws.Range("R23").Value = "ROC1"
$ws.Range("S23").Value = "2014-07-06T06:58:51"
$ws.Range("T23").Value = "2014-07-06T06:58:52"
# Row 24
$ws.Range("Q24").Value = "V25A"
$ws.Range("R24").Value = "CO02"
$ws.Range("S24").Value = "2014-07-16T03:38:01"
$ws.Range("T24").Value = "2014-07-16T03:38:03"
# Row 25
$ws.Range("Q25").Value = "AC05"
$ws.Range("R25").Value = "GO04"
$ws.Range("S25").Value = "2014-07-19T13:01:06"
$ws.Range("T25").Value = "2014-07-19T13:01:08"
# Row 26
$ws.Range("Q26").Value = "PB11"
$ws.Range("R26").Value = "PX03"
$ws.Range("S26").Value = "2014-07-28T10:22:52"
$ws.Range("T26").Value = "2014-07-28T10:22:52"
# Row 27
$ws.Range("Q27").Value = "PX03"
$ws.Range("R27").Value = "PB11"
$ws.Range("S27").Value = "2014-08-05T17:13:21"
$ws.Range("T27").Value = "2014-08-05T17:13:22"
# Row 28
$ws.Range("Q28").Value = "TLL"
$ws.Range("R28").Value = "CO03"
$ws.Range("S28").Value = "2014-08-05T21:02:14"
$ws.Range("T28").Value = "2014-08-05T21:02:16"
# Row 29
$ws.Range("Q29").Value = "AC07"
$ws.Range("R29").Value = "CO10"
$ws.Range("S29").Value = "2014-09-13T20:58:04"
$ws.Range("T29").Value = "2014-09-13T20:58:10"
# Row 30
$ws.Range("Q30").Value = "PB11"
$ws.Range("R30").Value = "PB08"
$ws.Range("S30").Value = "2014-09-14T13:14:31"
$ws.Range("T30").Value = "2014-09-14T13:14:32"
# Row 31
$ws.Range("Q31").Value = "VA01"
$ws.Range("R31").Value = "VA03"
$ws.Range("S31").Value = "2014-10-17T13:31:07"
$ws.Range("T31").Value = "2014-10-17T13:31:07"
# Row 32
$ws.Range("Q32").Value = "CO04"
$ws.Range("R32").Value = "MT02"
$ws.Range("S32").Value = "2014-11-28T09:27:41"
$ws.Range("T32").Value = "2014-11-28T09:27:42"
# Row 33
$ws.Range("Q33").Value = "AC05"
$ws.Range("R33").Value = "GO04"
$ws.Range("S33").Value = "2015-01-05T22:09:33"
$ws.Range("T33").Value = "2015-01-05T22:09:36"
# Row 34
$ws.Range("Q34").Value = "PB07"
$ws.Range("R34").Value = "PB15"
$ws.Range("S34").Value = "2015-01-28T10:08:42"
$ws.Range("T34").Value = "2015-01-28T10:08:45"
# Row 35
$ws.Range("Q35").Value = "CO03"
$ws.Range("R35").Value = "GO04"
$ws.Range("S35").Value = "2015-02-07T18:42:03"
$ws.Range("T35").Value = "2015-02-07T18:42:04"
# Row 36
$ws.Range("Q36").Value = "PATCX"
$ws.Range("R36").Value = "PB11"
$ws.Range("S36").Value = "2015-02-14T15:00:23"
$ws.Range("T36").Value = "2015-02-14T15:00:23"
# Row 37
$ws.Range("Q37").Value = "AC06"
$ws.Range("R37").Value = "GO03"
$ws.Range("S37").Value = "2015-03-01T07:46:45"
$ws.Range("T37").Value = "2015-03-01T07:46:49"
# Row 38
$ws.Range("Q38").Value = "TLL"
$ws.Range("R38").Value = "AC05"
$ws.Range("S38").Value = "2015-04-07T12:15:18"
$ws.Range("T38").Value = "2015-04-07T12:15:19"
# Row 39
$ws.Range("Q39").Value = "IN41"
$ws.Range("R39").Value = "GO04"
$ws.Range("S39").Value = "2015-08-22T08:24:52"
$ws.Range("T39").Value = "2015-08-22T08:24:56"
# Row 40
$ws.Range("Q40").Value = "PB11"
$ws.Range("R40").Value = "PSGCX"
$ws.Range("S40").Value = "2020-06-19T05:39:26"
$ws.Range("T40").Value = "2020-06-19T05:39:26"
